$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Fix the two "coneccion" -> "conexión" typos and drop the now
#    obsolete spell-check proofErr wrappers around them (the
#    w:proofErr "gramStart/gramEnd" pair around "para" is left alone).
#    Extending the replaced range one word past each proofErr boundary
#    forces the run containing the fixed word to merge with its
#    neighbour, which drops the proofErr markers that used to sit on
#    that boundary.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("por coneccion de", $false, $false, $false, $false, $false, $true, 1, $false, "por conexión de", 2)

$rng2 = $d.Content
$rng2.Find.Execute("la coneccion de", $false, $false, $false, $false, $false, $true, 1, $false, "la conexión de", 2)

# ---------------------------------------------------------------------
# 2. Drop the trailing "Escenario de Excepciones" row and the empty
#    placeholder row that followed it.
# ---------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Rows.Item($t.Rows.Count).Delete()
$t.Rows.Item($t.Rows.Count).Delete()

# ---------------------------------------------------------------------
# 3. Relocate the "_GoBack" bookmark from the title line to the end of
#    the "Escenario Básico Típico" table cell.
#
#    A genuinely collapsed Range positioned exactly at the end of a
#    run (with nothing following it inside the same paragraph) cannot
#    be anchored directly, so a one-character placeholder run is
#    inserted first to give the engine something to anchor the
#    collapsed bookmark range against; the placeholder is then removed
#    again, leaving the bookmark sitting cleanly right after the run.
# ---------------------------------------------------------------------
$cellRng = $d.Content
$cellRng.Find.Execute("Escenario Básico Típico")
$endPos = $cellRng.End

$anchor = $d.Range($endPos, $endPos)
$anchor.InsertAfter("X")

$bmRng = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

$placeholder = $d.Range($endPos, $endPos + 1)
$placeholder.Delete()
